$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notes")

# --- Update Description (row 2) ---
$ws.Range("A2").Value = "Description: Deprivation Of A Decent Standard Of Living"

# --- Update Source (row 4) ---
$ws.Range("A4").Value = "Source: Profiles of higher local governments 2014 - Uganda Bureau of Statistics"

# --- Insert a new row for the Source-link, right after Source (row 4) ---
$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = "Source-link: http://www.ubos.org/onlinefiles/uploads/ubos/2009_HLG_%20Abstract_printed/CIS+UPLOADS/Profiles%20of%20Higher%20Local%20Governments_June_2014.pdf"

# --- Update the license line (now at row 14 after the insert above) ---
$ws.Range("A14").Value = "It is licensed under a Creative Commons Attribution 4.0 International license."

# --- Insert a new row for the licensing link, right after that line (row 15) ---
$ws.Rows.Item(15).Insert()
$ws.Range("A15").Value = "More information on licensing is available here: https://creativecommons.org/licenses/by/4.0/"
